$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.630.45'
$ws.Range('E2').Value = '  -1.02%  '
$ws.Range('D3').Value = '3.782.47'
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.61'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.77'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.81%  '
$ws.Range('D7').Value = '3.781.22'
$ws.Range('E7').Value = '  +1.03%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('E10').Value = '  +0.43%  '
$ws.Range('E11').Value = '  -1.54%  '
$ws.Range('E12').Value = '  +0.36%  '
$ws.Range('E13').Value = '  -1.93%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.94'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.36%  '
$ws.Range('D15').Value = '4.416.78'
$ws.Range('E15').Value = '  +1.01%  '
$ws.Range('D16').Value = '3.820.75'
$ws.Range('E16').Value = '  +1.80%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.46'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +3.21%  '
$ws.Range('D18').Value = '67.621.22'
$ws.Range('E18').Value = '  -1.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.02'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.06%  '
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('E21').Value = '  -6.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '459.55'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.62%  '
$ws.Range('E23').Value = '  +0.49%  '
$ws.Range('E24').Value = '  +5.61%  '
$ws.Range('E25').Value = '  -0.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.03'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.73%  '
$ws.Range('E27').Value = '  -1.14%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.02'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('B29').Value = 'Dai'
$ws.Range('C29').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('E30').Value = '  +0.51%  '
$ws.Range('E31').Value = '  +3.85%  '
$ws.Range('E32').Value = '  -0.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '29.59'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.14%  '
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('E36').Value = '  +0.34%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.33'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('E38').Value = '  +0.57%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.994'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.38%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.75'
$ws.Range('D40').ClearFormats()
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '46.10'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +6.98%  '
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '48.19'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.57%  '
$ws.Range('E45').Value = '  -0.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '149.41'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +3.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.32'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '392.25'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.64%  '
$ws.Range('E49').Value = '  -4.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '26.47'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.97%  '
$ws.Range('D51').Value = '2.723.01'
$ws.Range('E51').Value = '  -1.16%  '
